# Updates translations and K&M UG Control button labels.
# Adds a new "DISCHARGE" label (using the same "\\f" substitution code as
# "FILL") to the Labels sheet, right below the existing FILL row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Labels")

# Insert a new row right after row 27 ("\\f" / "FILL"), pushing the
# remaining rows (RELEASE, HEATING, COOLING, FLAP, CONTROL, ...) down by one.
$ws.Rows("28:28").Insert()

$ws.Range("A28").Value = "\\f"
$ws.Range("B28").Value = "DISCHARGE"

# Match the author's recorded selection after the edit.
$ws.Range("B27").Select()
